$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1. Relocate the two trailing rows (totals row, footer row) two rows further
#    down, to make room for two new item rows. Work bottom-up so a target
#    range is never clobbered before it has been read from.
# ---------------------------------------------------------------------------

# old row 9 (footer: timestamp / page / credit) -> new row 11
$ws.Range("A9:Q11").UnMerge()
$ws.Range("A9:Q9").Copy()
$ws.Range("A11:Q11").PasteSpecial($xlPasteFormats)
$ws.Rows.Item(11).RowHeight = $ws.Rows.Item(9).RowHeight

# old row 8 (totals) -> new row 10
$ws.Range("A8:Q8").UnMerge()
$ws.Range("A8:Q8").Copy()
$ws.Range("A10:Q10").PasteSpecial($xlPasteFormats)
$ws.Rows.Item(10).RowHeight = $ws.Rows.Item(8).RowHeight

# ---------------------------------------------------------------------------
# 2. Turn rows 8 and 9 into item rows, matching row 7's layout/format.
# ---------------------------------------------------------------------------
$ws.Range("A7:Q7").Copy()
$ws.Range("A9:Q9").PasteSpecial($xlPasteFormats)
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q8").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 24.75
$ws.Rows.Item(9).RowHeight = 25.5

# ---------------------------------------------------------------------------
# 3. Re-create merges for rows 8, 9, 10, 11 (row 7's merges are untouched).
# ---------------------------------------------------------------------------
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

$ws.Range("P10:Q10").Merge()

$ws.Range("A11:F11").Merge()
$ws.Range("G11:I11").Merge()
$ws.Range("K11:Q11").Merge()

# ---------------------------------------------------------------------------
# 4. Write the new cell values.
# ---------------------------------------------------------------------------

# Row 7: new item - CARVID 6.25MG 30TAB
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "CARVID 6.25MG 30TAB"
$ws.Range("H7").Value = "1:0"
$ws.Range("L7").Value = "0"
$ws.Range("N7").Value = "45.00"
$ws.Range("P7").Value = "45.0000"
$ws.Range("Q7").Value = "1:0"

# Row 8: existing item (FAYCID), shifted down
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "FAYCID HAIR OIL 60 ML"
$ws.Range("H8").Value = "8:0"
$ws.Range("L8").Value = "1"
$ws.Range("N8").Value = "30.00"
$ws.Range("P8").Value = "30.0000"
$ws.Range("Q8").Value = "1:0"

# Row 9: new item - GLIMET FORTE 5/800 MG 30 F.C.TAB.
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "GLIMET FORTE 5/800 MG 30 F.C.TAB."
$ws.Range("H9").Value = "2:2"
$ws.Range("L9").Value = "1"
$ws.Range("N9").Value = "66.00"
$ws.Range("P9").Value = "66.0000"
$ws.Range("Q9").Value = "1:0"

# Row 10: totals
$ws.Range("P10").Value = 141

# Row 11: footer
$ws.Range("A11").Value = "Wednesday, 8 October, 2025 9:37 AM"
$ws.Range("G11").Value = "1/1"
$ws.Range("K11").Value = "developed by : Abdelaziz Talaat"
